$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): Location, other location, Pincode
$ws.Range("E1").Value = "Location"
$ws.Range("F1").Value = "other location"
$ws.Range("G1").Value = "Pincode"

# New data cells (row 2): Chennai, newyork, 1234567
# (write F2 first so the shared-string table picks up "newyork" before
# "Chennai", matching the target string order)
$ws.Range("F2").Value = "newyork"
$ws.Range("E2").Value = "Chennai"
$ws.Range("G2").Value = 1234567

# Resize existing column E and add width for new column F
# (ColumnWidth is quantised to whole pixels by the host, same as real Excel,
# so these land on the nearest representable width to the authored file's
# 12.54296875 / 13.26953125 character widths)
$ws.Columns.Item(5).ColumnWidth = 11.67
$ws.Columns.Item(6).ColumnWidth = 12.5

# Move the active selection to C10, matching the saved view state
$null = $ws.Range("C10").Select()
